$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.181.63'
$ws.Range("E2").Value = '  -2.13%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.010.76'
$ws.Range("E3").Value = '  -4.21%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.96'
$ws.Range("E5").Value = '  -3.15%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.77'
$ws.Range("E6").Value = '  -5.03%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.007.76'
$ws.Range("E8").Value = '  -4.15%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.135'
$ws.Range("E10").Value = '  -4.77%  '

$ws.Range("E11").Value = '  -1.21%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.430'
$ws.Range("E12").Value = '  -5.12%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000222'
$ws.Range("E13").Value = '  -4.91%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.80'
$ws.Range("E14").Value = '  -3.20%  '

$ws.Range("E15").Value = '  -0.20%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.508.57'
$ws.Range("E16").Value = '  -3.61%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.245.52'
$ws.Range("E17").Value = '  -2.01%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.011.77'
$ws.Range("E18").Value = '  -4.51%  '

$ws.Range("E19").Value = '  -5.00%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '437.84'
$ws.Range("E20").Value = '  -3.27%  '

$ws.Range("E21").Value = '  -5.67%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.661'
$ws.Range("E22").Value = '  -5.28%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.14'
$ws.Range("E23").Value = '  -5.65%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.93'
$ws.Range("E24").Value = '  -5.34%  '

$ws.Range("E25").Value = '  -5.93%  '

$ws.Range("E26").Value = '  -0.08%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.28%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.50'
$ws.Range("E28").Value = '  -6.36%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.18'
$ws.Range("E29").Value = '  -6.44%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.21'
$ws.Range("E30").Value = '  -7.34%  '

$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.56'
$ws.Range("E31").Value = '  -5.69%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.87'
$ws.Range("E32").Value = '  -6.61%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0938'
$ws.Range("E33").Value = '  -8.67%  '

$ws.Range("E34").Value = '  -4.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.951'
$ws.Range("E35").Value = '  -7.55%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.55'
$ws.Range("E36").Value = '  -3.87%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '50.12'
$ws.Range("E37").Value = '  -1.86%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0676'
$ws.Range("E38").Value = '  -3.44%  '

$ws.Range("E39").Value = '  -5.58%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.72'
$ws.Range("E40").Value = '  -3.61%  '

$ws.Range("E41").Value = '  -2.29%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '373.98'
$ws.Range("E42").Value = '  -5.57%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.645.25'
$ws.Range("E43").Value = '  -3.56%  '

$ws.Range("E44").Value = '  -9.16%  '

$ws.Range("E46").Value = '  -5.39%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '119.40'
$ws.Range("E47").Value = '  -4.40%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.97'
$ws.Range("E48").Value = '  -6.83%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '33.02'
$ws.Range("E49").Value = '  -3.76%  '

$ws.Range("E50").Value = '  -3.86%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.58'
$ws.Range("E51").Value = '  -6.21%  '
